$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data values for row 2
$ws.Range("G2").Value = 15710
$ws.Range("H2").Value = 21072
$ws.Range("J2").Value = 157.26

# Update the data values for row 3
$ws.Range("G3").Value = 15710
$ws.Range("H3").Value = 21072
$ws.Range("J3").Value = 157.26

# Update the selected cell/range shown in the sheet view
$ws.Range("K7").Select()
